# TC_35.xlsx: bring the local copy in line with the "dev" branch version.
#
# Note: the source revision also drops the orphaned customXml/item1.xml +
# customXml/itemProps1.xml parts (a stray CEIC "MetadataLink" part that
# isn't wired into _rels/.rels, so Excel's CustomXMLParts collection never
# sees it -- there is no object-model surface to reach/delete it from
# COM/VBA). That part of the change is left to packaging/CI, not to an
# in-app edit, so it's out of scope here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet from "My Series" to "Data"
$ws.Name = "Data"

# 2) Update cell A11 text from "Function Description" to "Function Information"
$ws.Range("A11").Value = "Function Information"

# 3) Update the custom number format (numFmtId 166) applied to B27:B43
#    from "0.000" to "###0.000"
$ws.Range("B27:B43").NumberFormat = "###0.000"

# 4) Update the cell-A1 comment payload (CEIC Data Manager metadata blob)
$comment = $ws.Range("A1").Comment
$newCommentText = "ER0AAB+LCAAAAAAAAAOlWVtvG8cV/isLPrVAyV2Ski0Low14k0KUFAWSqqy8FMPdkTj1XtidWUl8S4EEKdIURVE4RXpFn1IUqGu0CZDavfyXwJLdp/6FnrnsjaRirWoI1s4535k5c+bcZoTeufI944JEjIbBXqVesyoGCZzQpcH5XiXmZ9X6g8o7NupdOcQ7whH2CQewAVIB271idK8y53yxa5qXl5e1y2YtjM7NhmXVzcfDwcSZEx9XacA4DhxSSaXct0tVbNRx/SHh2MUcK8m9Sn/Sr3UIdbpAG+IAn5Oo1o4ZDQhjvYBTTgkTkhHBnHS6w++pjdmN2oNaHZlr9AzZjqnnKlwBqegaB8uSKfWJ3bDqO1Vrp9qsT62Hu8367tZObXt7671EMAWiAWZ8QqIL6kjChGN/IcWtnWbdetisb+0gcyMI5soMYKOR547JBWXE7RDPY6UsYuoDbDkcdl3OmBYyc7J6ovurcBDhxXxKuUfKqTEetg0/0Lpkk9hoP4yIA/a7l0qH5HIUabNOFwPgTuc04ssuXpae65iRaLQQRionaqNuGPCWRyJ+vICzJi64AjBsHsUEmbcwM6EuZQ580yAmrn2GPZYXKjDRSRg9YQvskEOIY1PMcRl4IXbB4ThlnDrZomsMdBSFC5gRFm+HnrsPs2rwBkY6cz8AE4tl22H4JNNuExPJU5XnC2fqY57A1+hoMg8vR4G3nMQz5kR0RtxuO0Fv5CERkFq6EzMe+qBFRkKKlqMs4R8E4CoZdYlDfewdeWBEZjdhlgIBtWIenlHeCb3YD1ii0woVncCOpuQq3WE6RiM43EAYPQz6QYJXZt7IKgqMw8t0zXWGNEKO3GJOctzrjFVwF2jJ8a1z5ImIXe5TD+pD/ixy1KJXTOaE8I0uoThIpMJ9UXHsluchMxsicEtwbVDDtqBoVOXP1LJ25Q8sm7JRL3BvxyVMdBj7oxmE74XckF0H3goJwRa8toeDJ0A9oXx+2EpU38BBasO34td5CAJ14eGlJKdGydNQP3C82CUq/vvBmfRIoZs6w1vZaI00gJC2EQ6W0+UC0jCjuxw+9ipQmHcZj6D0V2wnjAMeLUWiQKaGvk2GxbNALoC9O8ucReSHMXQcy/04cDqhe/fVXGWd44Dyu2sYxpHKfncXkdYTiTBmXSJSiszxd5Z3yuyJRaXgfkD8MKDO3a0NRhbau/fYCEui6s4SRMXXnfEelHFV5USs31ksgnYR6lqpZVqMhQ6VzqrDw83Jm7eETJec4diDVo1DRT1PU+0qGbXYk1VMnoSOIy9JeLZohBl0wo7r1xzoFUS3V3NCXxBMaEBPJsjM40XD45BecD7AwXkMLUWaV1bpaboV5XAa4YCJ7aQdxErm3QxCSZ5SnY29nuyQuQJBU+Ivwgh7Q7AJ3dcepzsjaDqGmM/1CKqYR5zEvmYmmkoVlUp0fhtMliO1AxHrOkOuECVI+KdqtzNMRkNii0OISK+DPTqLVEJNivYmHpxV1gYmqVdsrmRLmJgfrlxQZ79LlqIPzwaaLr21njCU74ocak/GWzuNbavZgB5GjJHc8Zhgz+hBHHNi9IMLwrgPYrvGmDDqwhfF3q7xLpkRCvVPmkhXoNLSeTm0n6R4qUoL9C1SigDoLM4pVJB1YMrJBOxTgiNvmQOqrQ5CB3A3P/739W+ev3rx+c3HT998+cF///GrV//82fWzD+Hj5q9/u/7kl2qbCoymeOYRqdC0vbNjNbfAz1ISEsY1ZRPsxg6XtNNT2fumY6SvcHLQ6fU7B4O2TCUpMRFX1cQUt8NlGGfDidqEXEgeqZl4goLY0yQ16XGBm6tOtritXZAiOs+/TVDZ4vXLz1+//POt0tpgWZtVf/Rou1pvvLULg3tvfQ2XdmGDQvoX4K2qtV1tNHLgFQwaq9yf2qnv2nCtfmQ1mlY9TeNu6sibQKssPdMUn5srcorUUZ1R6gL5ccKUjj+FEEnZKhRyA+2iX/7kzV+eFlDauppSnAWUkx2MWMxMBnLqw/HUmIyOx52eMe1NhJ9kvBxOTf4NYL16Gk8FpwqCGHvfMaCuQx0zKnD3qRjhmUGwMzeWEIm5OCw42yaqWuieU65qeRCF8UKdSE4go25Aptlko8SGXCN50p5rSSdjbYArXa///sUmAb2RbtbLpi8deRoqcBQpx9dR+9m/Xn310asXL26e//z6qx8VZtDrpA8A4OcQTflh6vaQ8nS9WaGgk4k05hPr+7n6ooniDnUU0oAzu/5AXp/0CIFoXcwmf6O+DyVPTiztBfQVCnoXs94V14FtHyKzSAA9FxiqbZjdMlOCyuGZXf/z29/d/PqLm0+fv/noT9cf//H6k09fv/z9m2d/UFF38/T5zU+f6Sy/WgikLuLuqvo/Q76EOIaIRkPUbuPr939hBCE3oOUwYpmRvn7/s9xkQlHZnGQzQzeXKlJUYQ2aFxZyRk6VVIeCXCqiGoCOKGHNFKGLWLigTrbIe1UxlYg7yfhWf1qNGTFC6Ka+DTspgjPhu8ppEVVSjx5ajXpDc5U2YgszzHKmP/DCGTQZCUM+NaxAClLfLJBh5XoHg1G7NcggSolR5JJIuKH6QElLKUpKnyWjxNVyFOBC4+fEnngdWoOts9KZc2nM1M8sZy1XpL/NrxQFBOrEUaQaokC/2k/iBTTDyWPc7Xz5QJnrfw9Vr5rviLNxv1vkwzjHhUJYZAuC5MvUpFkqTfWZeNFR7eyhME02BF7hURPMoR/mVad1AX1lZIq804uiMNqYfDJOAhtCJw0ZxcwsnmLkmaqu283OKiEkCS/9UJc+vcOwSzzCy71am5n0MLy4tyycfVnRPht5rjZmuatHapZsgvzTvXCU//flXjlbK4qgsRJPfaWf2pM76xiuuiW1UVuRguIGCKvrt/B9GjH+WGQC/aUopynlVHWoj8WFS33I8am91VAEAJj52c2CmknocvUHlNAbUJ+WvBZaSXwXJwFbLhaqheuX8xRRWg7JFTSYuRkgKc5+AGVDPaGUmU05LOTSVF48WzJ6PudlFXs4w8QlM6vqzEijuuVaO9VHhDSr9Tr8j51Gw7K2xaOnnhwyByWXJRcxkwPL/rBp/w8lieNLER0AAA=="
$comment.Text($newCommentText) | Out-Null
